$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.288846492767334
$ws.Range("B1").Value = 2.179458141326904
$ws.Range("C1").Value = 4.702958106994629
$ws.Range("D1").Value = 3.221054315567017
$ws.Range("E1").Value = 1.350842952728271
